# act tablas web jul25
# Adds 2023/2022 data points to the "Data" sheet (shifting historical rows down
# by two positions) and records an "actualizacion" / "Julio 2025" entry on the
# "Metadata" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("Metadata")

# ---------------------------------------------------------------------------
# Sheet "Data": insert two new rows right after the header row so every
# existing year moves down two rows, then stamp the final values for every
# row (years in column A as text, values in column B as numbers).
# ---------------------------------------------------------------------------
$ws1.Rows.Item(2).Insert()
$ws1.Rows.Item(2).Insert()

$dataYears = @("2023", "2022", "2021", "2019", "2018", "2017", "2016", "2015", "2014", "2013", "2012", "2011", "2010", "2009", "2008", "2007", "2006")
$dataValues = @(99, 99, 98.6, 98.8, 98.7, 98.8, 98.6, 98.6, 98.3, 97.8, 97.8, 96.5, 95.9, 95.4, 94.7, 94.4, 93.4)

$yearsRange = $ws1.Range("A2:A18")
$yearsRange.NumberFormat = "@"
for ($i = 0; $i -lt $dataYears.Length; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 1).Value = $dataYears[$i]
    $ws1.Cells.Item($row, 2).Value = $dataValues[$i]
}
$yearsRange.ClearFormats()

# ---------------------------------------------------------------------------
# Sheet "Metadata": the blank placeholder in A1 becomes a single space (to
# match B1), and a new "actualizacion" / "Julio 2025" row is inserted right
# before the "cita" row.
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = " "

$ws2.Rows.Item(9).Insert()
$ws2.Range("A9").Value = "actualizacion"
$ws2.Range("B9").Value = "Julio 2025"
